# aggiornamento fino a 9 agosto 2021
# Appends new daily rows (329-343) to the FinaleEmilia report, continuing the
# series through 2021-08-09, mirroring the formatting of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: row, date-serial (col A), nuovi pos. (col B), somma mobile 7gg. (col C),
# somma mobile 7gg. per 100mila abitanti (col D)
$data = @(
  @(329, 44403, 0, 3, 19.79414093428345),
  @(330, 44404, 1, 4, 26.39218791237794),
  @(331, 44405, 0, 4, 26.39218791237794),
  @(332, 44406, 0, 4, 26.39218791237794),
  @(333, 44407, 1, 2, 13.19609395618897),
  @(334, 44408, 0, 2, 13.19609395618897),
  @(335, 44409, 0, 2, 13.19609395618897),
  @(336, 44410, 0, 2, 13.19609395618897),
  @(337, 44411, 0, 1, 6.598046978094485),
  @(338, 44412, 0, 1, 6.598046978094485),
  @(339, 44413, 1, 2, 13.19609395618897),
  @(340, 44414, 1, 2, 13.19609395618897),
  @(341, 44415, 0, 2, 13.19609395618897),
  @(342, 44416, 2, 4, 26.39218791237794),
  @(343, 44417, 0, 4, 26.39218791237794)
)

# Column A on the existing data carries a date-style format (s="2" / centered,
# bordered, date numFmt). Replicate that style onto the new A329:A343 cells by
# copying the formatting from the last existing date cell (A328).
$ws.Range("A328").Copy()
$ws.Range("A329:A343").PasteSpecial(-4122)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
